$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format to avoid Excel auto-converting
# numeric-looking strings (e.g. "1.003") into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.593.94"
$ws.Range("E2").Value = "  -7.42%  "

$ws.Range("D3").Value = "1.690.05"
$ws.Range("E3").Value = "  -6.38%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "218.52"
$ws.Range("E5").Value = "  -5.86%  "

$ws.Range("D6").Value = "0.5053"
$ws.Range("E6").Value = "  -14.83%  "

$ws.Range("D8").Value = "0.2636"
$ws.Range("E8").Value = "  -5.23%  "

$ws.Range("D9").Value = "22.06"
$ws.Range("E9").Value = "  -5.58%  "

$ws.Range("D10").Value = "0.06227"
$ws.Range("E10").Value = "  -8.99%  "

$ws.Range("D11").Value = "0.07297"
$ws.Range("E11").Value = "  -2.53%  "

$ws.Range("D12").Value = "1.687.72"
$ws.Range("E12").Value = "  -6.70%  "

$ws.Range("D13").Value = "4.467"
$ws.Range("E13").Value = "  -6.80%  "

$ws.Range("D14").Value = "0.5800"
$ws.Range("E14").Value = "  -7.11%  "

$ws.Range("D15").Value = "1.915.96"
$ws.Range("E15").Value = "  -6.60%  "

$ws.Range("D16").Value = "0.000008242"
$ws.Range("E16").Value = "  -11.75%  "

$ws.Range("D17").Value = "65.19"
$ws.Range("E17").Value = "  -14.04%  "

$ws.Range("D18").Value = "26.629.46"
$ws.Range("E18").Value = "  -7.11%  "

$ws.Range("D19").Value = "5.033"
$ws.Range("E19").Value = "  -8.43%  "

$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").Value = "10.82"
$ws.Range("E21").Value = "  -5.74%  "

$ws.Range("D22").Value = "186.13"
$ws.Range("E22").Value = "  -11.97%  "

$ws.Range("D23").Value = "6.222"
$ws.Range("E23").Value = "  -9.41%  "

$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "144.80"
$ws.Range("E25").Value = "  -6.16%  "

$ws.Range("D26").Value = "7.546"
$ws.Range("E26").Value = "  -4.32%  "

$ws.Range("D27").Value = "0.1140"
$ws.Range("E27").Value = "  -10.52%  "

$ws.Range("D28").Value = "15.57"
$ws.Range("E28").Value = "  -5.34%  "

$ws.Range("D29").Value = "1.302"
$ws.Range("E29").Value = "  -8.40%  "

$ws.Range("D30").Value = "0.05699"
$ws.Range("E30").Value = "  -8.63%  "

$ws.Range("D31").Value = "1.333"
$ws.Range("E31").Value = "  -6.38%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "3.503"
$ws.Range("E32").Value = "  -6.87%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "3.493"
$ws.Range("E33").Value = "  -7.72%  "

$ws.Range("D34").Value = "1.642"
$ws.Range("E34").Value = "  -5.02%  "

$ws.Range("D35").Value = "1.014"
$ws.Range("E35").Value = "  -5.10%  "

$ws.Range("D36").Value = "0.5944"
$ws.Range("E36").Value = "  -7.38%  "

$ws.Range("D37").Value = "2.371"
$ws.Range("E37").Value = "  -4.85%  "

$ws.Range("D38").Value = "2.678"
$ws.Range("E38").Value = "  -1.39%  "

$ws.Range("D39").Value = "0.01596"
$ws.Range("E39").Value = "  -7.15%  "

$ws.Range("D40").Value = "1.074.21"
$ws.Range("E40").Value = "  -5.89%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.8657"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.910"
$ws.Range("E42").Value = "  -9.24%  "

$ws.Range("E43").Value = "  -0.80%  "

$ws.Range("D44").Value = "98.18"
$ws.Range("E44").Value = "  -2.19%  "

$ws.Range("D45").Value = "1.843.93"
$ws.Range("E45").Value = "  -6.02%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000107"
$ws.Range("E46").Value = "  -4.36%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "56.53"
$ws.Range("E47").Value = "  -6.68%  "

$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.49%  "

$ws.Range("D49").Value = "8.082"
$ws.Range("E49").Value = "  -3.77%  "

$ws.Range("D50").Value = "0.4314"
$ws.Range("E50").Value = "  -3.79%  "

$ws.Range("D51").Value = "0.05209"
$ws.Range("E51").Value = "  -4.86%  "

# Restore the default (General/Normal) style on column D so that cells
# keep matching their original (unstyled) appearance while remaining text.
$ws.Range("D2:D51").Style = "Normal"
